$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Complete row 10: add the missing X10 / Y10 values ---
$ws.Cells.Item(10, 24).Value = -3.4100040000000149   # X10
$ws.Cells.Item(10, 25).Value = "Down"                 # Y10 (shared string index 28)

# --- Append new row 11 ---
$ws.Cells.Item(11, 1).Value  = 42654.886666666665     # A11 (Date)
$ws.Cells.Item(11, 2).Value  = 12                      # B11
$ws.Cells.Item(11, 3).Value  = "Buy"                   # C11 (shared string index 25)
$ws.Cells.Item(11, 4).Value  = 34                      # D11
$ws.Cells.Item(11, 5).Value  = 20483                   # E11
$ws.Cells.Item(11, 6).Value  = 1171                    # F11
$ws.Cells.Item(11, 7).Value  = 59                      # G11
$ws.Cells.Item(11, 8).Value  = 39                      # H11
$ws.Cells.Item(11, 9).Value  = 94                      # I11
$ws.Cells.Item(11, 10).Value = 4                       # J11
$ws.Cells.Item(11, 11).Value = 35576                   # K11
$ws.Cells.Item(11, 12).Value = 156                     # L11
$ws.Cells.Item(11, 13).Value = 105                     # M11
$ws.Cells.Item(11, 14).Value = 60                      # N11
$ws.Cells.Item(11, 15).Value = 3                       # O11
$ws.Cells.Item(11, 16).Value = "Named"                 # P11 (shared string index 26)
$ws.Cells.Item(11, 17).Value = 39.313912976930268      # Q11
$ws.Cells.Item(11, 18).Value = 1.8                     # R11
$ws.Cells.Item(11, 19).Value = 0.0864                  # S11
$ws.Cells.Item(11, 20).Value = -0.0115                 # T11
$ws.Cells.Item(11, 21).Value = 5.85                    # U11
$ws.Cells.Item(11, 22).Value = "N/A"                   # V11 (shared string index 27)
$ws.Cells.Item(11, 23).Value = 0                       # W11

# Match number formats used by the same columns in the rows above (date / percent)
$ws.Cells.Item(11, 1).NumberFormat  = "m/d/yy h:mm"
$ws.Cells.Item(11, 19).NumberFormat = "0.00%"
$ws.Cells.Item(11, 20).NumberFormat = "0.00%"
